# Commit: "Correct typos in simple voting results for user-based graphs"
#
# The user_simple_vot_25 block (results!D314:I337) had copy/paste typos in
# its P@k / R@k / RPrec-ish columns (G, H, I). Fix the values and widen the
# named range to match the (now wider) table, matching the author's edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("results")

# --- 1. Correct the mistaken values in columns G:I for rows 314-337 -------
$fixedValues = @{
    314 = @{ G = 0.96808510638297796; H = 0.76418439716312003; I = 0.48958333333333298 }
    315 = @{ G = 0.96808510638297796; H = 0.76418439716312003; I = 0.48958333333333298 }
    316 = @{ G = 0.96808510638297796; H = 0.76418439716312003; I = 0.47451241134751798 }
    317 = @{ G = 0.96808510638297796; H = 0.76418439716312003; I = 0.48958333333333298 }
    318 = @{ G = 0.96808510638297796; H = 0.76418439716312003; I = 0.48958333333333298 }
    319 = @{ G = 0.96808510638297796; H = 0.76418439716312003; I = 0.48958333333333298 }
    320 = @{ G = 0.96808510638297796; H = 0.76418439716312003; I = 0.48958333333333298 }
    321 = @{ G = 0.96808510638297796; H = 0.76418439716312003; I = 0.48958333333333298 }
    322 = @{ G = 0.98138297872340396; H = 0.69095744680851101; I = 0.48958333333333298 }
    323 = @{ G = 0.98138297872340396; H = 0.69095744680851101; I = 0.48958333333333298 }
    324 = @{ G = 0.98138297872340396; H = 0.69095744680851101; I = 0.47451241134751798 }
    325 = @{ G = 0.98138297872340396; H = 0.69095744680851101; I = 0.48958333333333298 }
    326 = @{ G = 0.98138297872340396; H = 0.69095744680851101; I = 0.48958333333333298 }
    327 = @{ G = 0.98138297872340396; H = 0.69095744680851101; I = 0.48958333333333298 }
    328 = @{ G = 0.98138297872340396; H = 0.69095744680851101; I = 0.48958333333333298 }
    329 = @{ G = 0.98138297872340396; H = 0.69095744680851101; I = 0.48958333333333298 }
    330 = @{ G = 0.99468085106382897; H = 0.58803191489361695; I = 0.48958333333333298 }
    331 = @{ G = 0.99468085106382897; H = 0.58803191489361695; I = 0.48958333333333298 }
    332 = @{ G = 0.99468085106382897; H = 0.59015957446808498; I = 0.47451241134751798 }
    333 = @{ G = 0.99468085106382897; H = 0.58803191489361695; I = 0.48958333333333298 }
    334 = @{ G = 0.99468085106382897; H = 0.58803191489361695; I = 0.48958333333333298 }
    335 = @{ G = 0.99468085106382897; H = 0.58803191489361695; I = 0.48958333333333298 }
    336 = @{ G = 0.99468085106382897; H = 0.58803191489361695; I = 0.48958333333333298 }
    337 = @{ G = 0.99468085106382897; H = 0.58803191489361695; I = 0.48958333333333298 }
}

foreach ($row in $fixedValues.Keys) {
    $cols = $fixedValues[$row]
    foreach ($col in $cols.Keys) {
        $ws.Range("$col$row").Value = $cols[$col]
    }
}

# --- 2. The corrected table is now wider; extend the named range ---------
#        results!user_simple_vot_25  $D$314:$I$337  ->  $D$314:$N$337
$wb.Names.Item("results!user_simple_vot_25").RefersTo = "=results!`$D`$314:`$N`$337"

# --- 3. Leave the selection where the author left it after the fix -------
$ws.Range("M332").Select()
